# "Added Indian MF 1st Stab" - insert 9 new weekly rating-history columns
# (Jun_16 .. Sep_08) immediately to the right of the report-name column (A),
# pushing the existing Jun_09..Mar_10 columns from B:V to K:AE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 9 new blank columns at B:J - this shifts the existing B:V block
#    (Jun_09 .. Mar_10) to K:AE, carrying values/styles/column widths along.
$ws.Columns("B:J").Insert()

# 2) New column headers for the inserted date columns (row 1), newest-first.
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# 3) The new columns default to "UN" (unchanged rating) for every analyst row,
#    same as the rest of the table.
$ws.Range("B2:J33").Value = "UN"

# 4) One real rating action landed in the new Jul_17 column for row 10
#    (Piper Jaffray Companies).
$ws.Range("F10").Value = "7/17/2019,Lowers Target,Overweight,$12.00 -> $10.00"
